# Insert a new weekly price-record row for "Perejil" (Vega Central Mapocho de
# Santiago) ahead of the existing row 447, pushing the remaining records
# (old rows 447:499) down to 448:500 -- matching a normal "new week's data
# appended near the top" edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 447:499 down to 448:500 and free up row 447 for the new record.
$ws.Rows(447).Insert()

# Populate the newly freed row 447 with the new weekly record.
$ws.Range("A447").Value = 9
$ws.Range("B447").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C447").Value = "Metropolitana"
$ws.Range("D447").Value = 44946
$ws.Range("E447").Value = 13
$ws.Range("F447").Value = 100112044
$ws.Range("G447").Value = "Perejil"
$ws.Range("H447").Value = "Sin especificar"
$ws.Range("I447").Value = "Primera"
$ws.Range("J447").Value = 70
$ws.Range("K447").Value = 15000
$ws.Range("L447").Value = 16000
$ws.Range("M447").Value = 15500
$ws.Range("N447").Value = "`$/docena de atados"
$ws.Range("O447").Value = "Región Metropolitana"
$ws.Range("P447").Value = 5167
$ws.Range("Q447").Value = 3
$ws.Range("R447").Value = "Hortaliza"
